$wb = $excel.ActiveWorkbook

# Sheet "展览" - update 想去人数 (attendance count) values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 389
$wsExhibit.Range("F4").Value = 3046
$wsExhibit.Range("F6").Value = 634

# Sheet "全部类型" - same rows duplicated, update matching values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 389
$wsAll.Range("F5").Value = 3046
$wsAll.Range("F7").Value = 634
